$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 15:52"

# Country name shifts caused by the country list being re-sorted by total cases
$ws.Range("A116").Value = "Republica de Yibuti"   # was "Sri Lanka"
$ws.Range("A117").Value = "Sri Lanka"   # was "Kenia"
$ws.Range("A118").Value = "Kenia"   # was "Mayotte"
$ws.Range("A119").Value = "Mayotte"   # was "Republica de Yibuti"
$ws.Range("A174").Value = "Republica del Chad"   # was "Islas Virgenes de los Estados Unidos"
$ws.Range("A175").Value = "Islas Virgenes de los Estados Unidos"   # was "Fiyi"
$ws.Range("A176").Value = "Fiyi"   # was "Namibia"
$ws.Range("A177").Value = "Namibia"   # was "Mongolia"
$ws.Range("A178").Value = "Mongolia"   # was "Dominica"
$ws.Range("A179").Value = "Dominica"   # was "Santa Lucia"
$ws.Range("A180").Value = "Santa Lucia"   # was "Granada"
$ws.Range("A181").Value = "Granada"   # was "Zimbabue"
$ws.Range("A182").Value = "Zimbabue"   # was "Suazilandia"
$ws.Range("A183").Value = "Suazilandia"   # was "Curazao"
$ws.Range("A184").Value = "Curazao"   # was "Botsuana"
$ws.Range("A185").Value = "Botsuana"   # was "Belice"
$ws.Range("A186").Value = "Malaui"   # was "San Cristobal y Nieves"
$ws.Range("A187").Value = "Belice"   # was "Nepal"
$ws.Range("A188").Value = "San Cristobal y Nieves"   # was "San Vicente y las Granadinas"
$ws.Range("A189").Value = "Nepal"   # was "Malaui"
$ws.Range("A190").Value = "San Vicente y las Granadinas"   # was "Seychelles"
$ws.Range("A191").Value = "Seychelles"   # was "Republica del Chad"
$ws.Range("A200").Value = "Cabo Verde"   # was "Santa Sede"
$ws.Range("A201").Value = "Santa Sede"   # was "Cabo Verde"

# Updated case numbers (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 533470
$ws.Range("C4").Value = 591
$ws.Range("D4").Value = 30523
$ws.Range("E4").Value = 482352
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 20595

$ws.Range("E30").Value = 6303
$ws.Range("F30").Value = 59
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 124

$ws.Range("B59").Value = 1701
$ws.Range("C59").Value = 12
$ws.Range("D59").Value = 889
$ws.Range("E59").Value = 804
$ws.Range("F59").Value = 10

$ws.Range("D62").Value = 94
$ws.Range("E62").Value = 1435
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 31

$ws.Range("D75").Value = 99
$ws.Range("E75").Value = 802

$ws.Range("D76").Value = 65
$ws.Range("E76").Value = 764

$ws.Range("E110").Value = 259
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 3

$ws.Range("B116").Value = 214
$ws.Range("C116").Value = 27
$ws.Range("D116").Value = 36
$ws.Range("E116").Value = 176
$ws.Range("F116").Value = 0
$ws.Range("H116").Value = 2

$ws.Range("B117").Value = 203
$ws.Range("C117").Value = 5
$ws.Range("D117").Value = 55
$ws.Range("E117").Value = 141
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 7

$ws.Range("B118").Value = 197
$ws.Range("C118").Value = 6
$ws.Range("D118").Value = 25
$ws.Range("E118").Value = 164
$ws.Range("F118").Value = 2
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 8

$ws.Range("B119").Value = 196
$ws.Range("D119").Value = 59
$ws.Range("E119").Value = 134
$ws.Range("F119").Value = 3
$ws.Range("H119").Value = 3

$ws.Range("B174").Value = 18
$ws.Range("C174").Value = 7
$ws.Range("D174").Value = 2
$ws.Range("E174").Value = 16

$ws.Range("B175").Value = 17
$ws.Range("E175").Value = 17

$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 16

$ws.Range("D177").Value = 3
$ws.Range("E177").Value = 13

$ws.Range("D178").Value = 4
$ws.Range("E178").Value = 12

$ws.Range("B179").Value = 16
$ws.Range("D179").Value = 5

$ws.Range("B180").Value = 15
$ws.Range("D180").Value = 4
$ws.Range("E180").Value = 11
$ws.Range("F180").Value = 0

$ws.Range("E181").Value = 14
$ws.Range("F181").Value = 2
$ws.Range("H181").Value = 0

$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 11
$ws.Range("H182").Value = 3

$ws.Range("C183").Value = 2
$ws.Range("E183").Value = 7
$ws.Range("H183").Value = 0

$ws.Range("B184").Value = 14
$ws.Range("D184").Value = 7
$ws.Range("E184").Value = 6

$ws.Range("E185").Value = 12
$ws.Range("F185").Value = 0
$ws.Range("H185").Value = 1

$ws.Range("B186").Value = 13
$ws.Range("C186").Value = 1
$ws.Range("E186").Value = 11
$ws.Range("F186").Value = 1
$ws.Range("H186").Value = 2

$ws.Range("B187").Value = 13
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("F187").Value = 1
$ws.Range("H187").Value = 2

$ws.Range("D188").Value = 0
$ws.Range("E188").Value = 12

$ws.Range("C189").Value = 3
$ws.Range("D189").Value = 1
$ws.Range("E189").Value = 11
$ws.Range("F189").Value = 0
$ws.Range("H189").Value = 0

$ws.Range("B190").Value = 12
$ws.Range("D190").Value = 1

$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 11

$ws.Range("D200").Value = 1
$ws.Range("H200").Value = 1

$ws.Range("D201").Value = 1
$ws.Range("H201").Value = 0
